$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Location) - header then data, top to bottom
$ws.Range("B1").Value = "Location"
$ws.Range("B2").Value = "Chennai"
$ws.Range("B3").Value = "Pune"
$ws.Range("B4").Value = "Delhi"

# Column C (Car Name) - header then data, top to bottom
$ws.Range("C1").Value = "Car Name"
$ws.Range("C2").Value = "Hyundai I10"
$ws.Range("C3").Value = "Maruti Swift"
$ws.Range("C4").Value = "Hyundai Santro Xing"

# Column D data first (no header yet)
$ws.Range("D2").Value = "hhasj"
$ws.Range("D3").Value = "S@n.com"
$ws.Range("D4").Value = "ma.com"

# Column D header added last
$ws.Range("D1").Value = "Email"

# Column widths
$ws.Columns.Item(3).ColumnWidth = 18
$ws.Columns.Item(4).ColumnWidth = 17

# Hyperlink on D3
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:S@n.com")

# Selection
[void]$ws.Range("F6").Select()
